# Optuna Attempt (go back with original)
# Update forecast values on the "Forecast Comparison" sheet and the
# corresponding roll-up metrics on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------
# Column D = MyForecast, H = Inventory Coverage, L = Seasonality Index

$wsForecast.Range("D2").Value = 299
$wsForecast.Range("H2").Value = 2.74
$wsForecast.Range("L2").Value = 0.89

$wsForecast.Range("D3").Value = 314
$wsForecast.Range("H3").Value = 1.66
$wsForecast.Range("L3").Value = 0.9

$wsForecast.Range("D4").Value = 328
$wsForecast.Range("H4").Value = 0.63
$wsForecast.Range("L4").Value = 0.89

$wsForecast.Range("D5").Value = 343
$wsForecast.Range("L5").Value = 1.11

$wsForecast.Range("D6").Value = 341
$wsForecast.Range("L6").Value = 1.08

$wsForecast.Range("L7").Value = 1.09
$wsForecast.Range("L8").Value = 0.99
$wsForecast.Range("L9").Value = 1.05
$wsForecast.Range("L10").Value = 1.03
$wsForecast.Range("L11").Value = 0.98
$wsForecast.Range("L12").Value = 1.08
$wsForecast.Range("L13").Value = 1.18
$wsForecast.Range("L14").Value = 1.16
$wsForecast.Range("L15").Value = 1.03
$wsForecast.Range("L16").Value = 0.92
$wsForecast.Range("L17").Value = 0.92

# --- Summary sheet ---------------------------------------------------------
# Values on this sheet are stored as text, so force the Text number format
# before assigning to avoid Excel auto-converting the strings to numbers.

$wsSummary.Range("B9:B11").NumberFormat = "@"
$wsSummary.Range("B14").NumberFormat = "@"

$wsSummary.Range("B9").Value = "5070"
$wsSummary.Range("B10").Value = "2619"
$wsSummary.Range("B11").Value = "1284"
$wsSummary.Range("B14").Value = "287"
